$d = $word.ActiveDocument

function Insert-RunsXmlAt($rangeStart, $innerXml) {
    $ins = $d.Range($rangeStart, $rangeStart)
    $wrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $ins.InsertXML($wrapped)
}

# The three m:xxx fields (written as classic Word fldChar/instrText fields) become
# plain "{...}" literal-text runs, as the template's field rewriter now emits
# the M2Doc tag delimiters directly instead of relying on Word field codes.

# Field 1: {m:if self.name = 'anydsl'}
$f = $d.Fields.Item(1)
$start = $f.Code.Start - 1
$f.Delete()
$para2 = '<w:r><w:t xml:space="preserve">{m:if </w:t></w:r><w:r><w:t xml:space="preserve">self.name </w:t></w:r><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>''</w:t></w:r><w:r><w:t>anydsl</w:t></w:r><w:r><w:t>''}</w:t></w:r>'
Insert-RunsXmlAt $start $para2

# Field 2: {m:'dh1.gif'.asImage().setWidth(100)} (keeps the bookmark + lang rPr on the 2nd run)
$f = $d.Fields.Item(1)
$start = $f.Code.Start - 1
$f.Delete()
$para3 = '<w:r><w:t>{m:''dh1.gif''.asImage()</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.setWidth(100)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">}</w:t></w:r>'
Insert-RunsXmlAt $start $para3

# Field 3: {m:endif}
$f = $d.Fields.Item(1)
$start = $f.Code.Start - 1
$f.Delete()
$para4 = '<w:r><w:t xml:space="preserve">{m:endif}</w:t></w:r>'
Insert-RunsXmlAt $start $para4

Write-Output "done"
